$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.2049062049062049
$ws.Range("C2").Value = 0.5425685425685426
$ws.Range("J2").Value = 0.01587301587301587
$ws.Range("P2").Value = 0.1341991341991342
$ws.Range("S2").Value = 0.1024531024531024

# Row 3
$ws.Range("B3").Value = 0.01036269430051814
$ws.Range("C3").Value = 0.0155440414507772
$ws.Range("J3").Value = 0.04404145077720207
$ws.Range("P3").Value = 0.7020725388601037
$ws.Range("S3").Value = 0.227979274611399

# Row 4
$ws.Range("J4").Value = 0.02912621359223301
$ws.Range("P4").Value = 0.6893203883495146
$ws.Range("S4").Value = 0.2815533980582524

# Row 6
$ws.Range("B6").Value = 0.0558252427184466
$ws.Range("D6").Value = 0.01213592233009709
$ws.Range("E6").Value = 0.002427184466019417
$ws.Range("F6").Value = 0.04854368932038835
$ws.Range("J6").Value = 0.2888349514563107
$ws.Range("O6").Value = 0.02669902912621359
$ws.Range("Q6").Value = 0.1820388349514563
$ws.Range("R6").Value = 0.07281553398058252
$ws.Range("S6").Value = 0.3106796116504854

# Row 7
$ws.Range("B7").Value = 0.115625
$ws.Range("D7").Value = 0.028125
$ws.Range("F7").Value = 0.053125
$ws.Range("J7").Value = 0.15625
$ws.Range("O7").Value = 0.009375
$ws.Range("Q7").Value = 0.19375
$ws.Range("R7").Value = 0.053125
$ws.Range("S7").Value = 0.390625

# Row 8
$ws.Range("B8").Value = 0.08073196986006459
$ws.Range("D8").Value = 0.01506996770721206
$ws.Range("F8").Value = 0.05920344456404737
$ws.Range("J8").Value = 0.1033369214208827
$ws.Range("O8").Value = 0.00968783638320775
$ws.Range("Q8").Value = 0.2055974165769645
$ws.Range("R8").Value = 0.09149623250807319
$ws.Range("S8").Value = 0.4348762109795479

# Row 9
$ws.Range("B9").Value = 0.1113744075829384
$ws.Range("D9").Value = 0.004739336492890996
$ws.Range("E9").Value = 0.002369668246445498
$ws.Range("F9").Value = 0.07109004739336493
$ws.Range("J9").Value = 0.1066350710900474
$ws.Range("O9").Value = 0.02369668246445497
$ws.Range("Q9").Value = 0.2061611374407583
$ws.Range("R9").Value = 0.07582938388625593
$ws.Range("S9").Value = 0.3981042654028436

# Row 10
$ws.Range("B10").Value = 0.1309949892627058
$ws.Range("D10").Value = 0.02755905511811024
$ws.Range("F10").Value = 0.05619183965640658
$ws.Range("J10").Value = 0.117394416607015
$ws.Range("O10").Value = 0.0161059413027917
$ws.Range("Q10").Value = 0.2269148174659986
$ws.Range("R10").Value = 0.08052970651395848
$ws.Range("S10").Value = 0.3443092340730136

# Row 11
$ws.Range("G11").Value = 0.1411992263056093
$ws.Range("J11").Value = 0.1063829787234043
$ws.Range("K11").Value = 0.1876208897485493
$ws.Range("L11").Value = 0.5570599613152805
$ws.Range("S11").Value = 0.007736943907156673

# Row 12
$ws.Range("G12").Value = 0.6677524429967426
$ws.Range("J12").Value = 0.241042345276873
$ws.Range("K12").Value = 0.01628664495114007
$ws.Range("L12").Value = 0.03908794788273615
$ws.Range("S12").Value = 0.03583061889250815

# Row 13
$ws.Range("F13").Value = 0.02531645569620253
$ws.Range("G13").Value = 0.6075949367088608
$ws.Range("J13").Value = 0.2911392405063291
$ws.Range("S13").Value = 0.0759493670886076

# Row 15
$ws.Range("F15").Value = 0.0148936170212766
$ws.Range("H15").Value = 0.1574468085106383
$ws.Range("I15").Value = 0.06808510638297872
$ws.Range("J15").Value = 0.3851063829787234
$ws.Range("K15").Value = 0.04680851063829787
$ws.Range("M15").Value = 0.00425531914893617
$ws.Range("O15").Value = 0.09148936170212765
$ws.Range("S15").Value = 0.2319148936170213

# Row 16
$ws.Range("F16").Value = 0.02558139534883721
$ws.Range("H16").Value = 0.1581395348837209
$ws.Range("I16").Value = 0.08372093023255814
$ws.Range("J16").Value = 0.4372093023255814
$ws.Range("K16").Value = 0.1093023255813954
$ws.Range("M16").Value = 0.02093023255813953
$ws.Range("O16").Value = 0.06046511627906977
$ws.Range("S16").Value = 0.1046511627906977

# Row 17
$ws.Range("F17").Value = 0.01821668264621285
$ws.Range("H17").Value = 0.1658676893576222
$ws.Range("I17").Value = 0.1045062320230105
$ws.Range("J17").Value = 0.4343240651965484
$ws.Range("K17").Value = 0.07861936720997123
$ws.Range("M17").Value = 0.01629913710450623
$ws.Range("O17").Value = 0.07094918504314478
$ws.Range("S17").Value = 0.1112176414189837

# Row 18
$ws.Range("F18").Value = 0.01288659793814433
$ws.Range("H18").Value = 0.2061855670103093
$ws.Range("I18").Value = 0.08247422680412371
$ws.Range("J18").Value = 0.4407216494845361
$ws.Range("K18").Value = 0.07731958762886598
$ws.Range("M18").Value = 0.01288659793814433
$ws.Range("O18").Value = 0.06958762886597938
$ws.Range("S18").Value = 0.09793814432989691

# Row 19
$ws.Range("F19").Value = 0.01215686274509804
$ws.Range("H19").Value = 0.2109803921568627
$ws.Range("I19").Value = 0.08313725490196078
$ws.Range("J19").Value = 0.4050980392156863
$ws.Range("K19").Value = 0.08941176470588236
$ws.Range("M19").Value = 0.01843137254901961
$ws.Range("N19").Value = 0.001568627450980392
$ws.Range("O19").Value = 0.06784313725490196
$ws.Range("S19").Value = 0.1113725490196078
